$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and two-row name/link swaps)

$ws.Range("D2").Value = "22.390.95"
$ws.Range("E2").Value = "  -4.77%  "
$ws.Range("D3").Value = "1.573.20"
$ws.Range("E3").Value = "  -4.61%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.08"
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("E7").Value = "  -3.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.61"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3369"
$ws.Range("E9").Value = "  -5.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07579"
$ws.Range("E11").Value = "  -6.53%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.12"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.055"
$ws.Range("E14").Value = "  -5.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.866"
$ws.Range("E15").Value = "  -7.59%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.583.48"
$ws.Range("E16").Value = "  -4.50%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001137"
$ws.Range("E17").Value = "  -5.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.46"
$ws.Range("E18").Value = "  -8.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06761"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.232"
$ws.Range("E21").Value = "  -8.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.32"
$ws.Range("E22").Value = "  -6.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.98"
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.426"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "22.413.45"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.966"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.79"
$ws.Range("E27").Value = "  -5.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "145.80"
$ws.Range("E28").Value = "  -4.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.929"
$ws.Range("E29").Value = "  -5.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.05"
$ws.Range("E30").Value = "  -6.22%  "
$ws.Range("D31").Value = "1.745.31"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.262"
$ws.Range("E32").Value = "  -9.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.981"
$ws.Range("E33").Value = "  -7.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9784"
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.40"
$ws.Range("E35").Value = "  -12.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08448"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02534"
$ws.Range("E37").Value = "  -7.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2299"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06521"
$ws.Range("E39").Value = "  -5.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.495"
$ws.Range("E40").Value = "  -8.43%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.260"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.75"
$ws.Range("E42").Value = "  -11.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6370"
$ws.Range("E43").Value = "  -8.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.42"
$ws.Range("E44").Value = "  -8.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5988"
$ws.Range("E46").Value = "  -7.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.777"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.118"
$ws.Range("E48").Value = "  -6.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.02"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07287"
$ws.Range("E50").Value = "  -6.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.185"
$ws.Range("E51").Value = "  -0.59%  "
